$wb = $excel.ActiveWorkbook

# Reorder sheets: "Ementas" should come before "caldendário e notas"
$wsCal = $wb.Worksheets.Item("caldendário e notas")
$wsEmentas = $wb.Worksheets.Item("Ementas")
$wsCal.Move($null, $wsEmentas)

# Re-fetch the grade sheet (worksheet refs are position-bound, so after the
# Move we need a fresh handle) and make it the active/selected tab.
$wsCal = $wb.Worksheets.Item("caldendário e notas")
$wsCal.Activate()

# Record exam grades for BIM 9D (row 13) and BIM 10D (row 14) in the
# "Prova" column (J). The dependent "Nota Final" formulas in column K
# recalculate automatically.
$wsCal.Range("J13").Value = 48
$wsCal.Range("J14").Value = 48

# Update the "Média geral" formula to also factor in the TCC row (15).
$wsCal.Range("K16").Formula = "=(K3+K4+K5+K6+K7+K8+K9+K10+K11+K12+K13+K14+K15)/13"
